$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# New row of data appended at the end of the table (row 13)
# (order matches the order new shared strings were interned in the source workbook)
# NOTE: these strings use a literal backslash+n (game-script line break marker),
# not an actual newline character - single-quoted so PowerShell does not expand it.
$ws.Range("C13").Value = ' [CS:P]Luminous Spring[CR]?[K] Oh, that place\nis where Pokémon used to go to evolve.'
$ws.Range("A13").Value = 'SCRIPT/G01P03A/us0104.ssb'
$ws.Range("D13").Value = ' [CS:P]Сияющий Источник[CR]?[K] О, это то\nместо, куда приходили Покемоны, чтобы\nэволюционировать.'
$ws.Range("E13").Value = ' [CS:P]Òéÿýþéê Éòóïœîéë[CR]?[K] Ï, üóï óï\níåòóï, ëôäà ðñéöïäéìé Ðïëåíïîú, œóïáú\nüâïìýøéïîéñïâàóû.'
$ws.Range("B13").Value = 57

# The new last row (13) takes on the plain formatting that row 12 had before the
# insertion, and row 12 takes on the "section separator" border formatting that
# row 11 (the previous last-of-section row) has.
$ws.Range("A12:E12").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("A13:E13").RowHeight = 43.2

$ws.Range("D12").Select()
